# Release v0.1.0-beta: Fix validation errors and update canonical URL
#
# Applies updates to the "Metadata" sheet (Property/Value pairs) and the
# "Include #0" sheet (Concept/Description pairs) of the ValueSet workbook.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version: 1.0.0 -> 0.1.0
$meta.Range("B3").Value = "0.1.0"

# Status: active -> draft
$meta.Range("B6").Value = "draft"

# Experimental: (blank) -> false
# Leading apostrophe forces this to be stored as literal text "false"
# rather than being auto-converted to a Boolean TRUE/FALSE cell.
$meta.Range("B7").Value = "'false"

# Date: updated publish timestamp
$meta.Range("B8").Value = "2025-12-26T14:13:58+00:00"

# Description: (blank) -> new description text
$meta.Range("B11").Value = "Value set for patient housing status"

# --- Include #0 sheet (expansion concepts) --------------------------------
$inc = $wb.Worksheets.Item("Include #0")

# Row 2: concept code + new description
# Leading apostrophe keeps the numeric-looking code stored as text.
$inc.Range("A2").Value = "'266935003"
$inc.Range("B2").Value = "Housing lack"

# Row 3: concept code + new description
$inc.Range("A3").Value = "'224224003"
$inc.Range("B3").Value = "Lives in staffed home"
